$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("initialize_fiscalYears")
$ws.Range("A2").Value = 0
